# Apply edit: add I0 and IF columns (I and J) with header cells and data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 (bold, centered, bordered) into I1 and J1
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2..85, columns I (I0) and J (IF)
$data = @(
    (4,5),
    (8,8),
    (7,7),
    (8,8),
    (7,8),
    (7,8),
    (8,8),
    (6,7),
    (7,7),
    (7,7),
    (8,8),
    (8,8),
    (7,7),
    (6,6),
    (6,6),
    (8,8),
    (8,8),
    (8,8),
    (8,8),
    (7,7),
    (9,9),
    (5,6),
    (8,8),
    (7,7),
    (6,6),
    (2,2),
    (6,7),
    (7,7),
    (8,8),
    (4,5),
    (7,7),
    (6,6),
    (7,8),
    (7,7),
    (9,9),
    (6,7),
    (5,5),
    (5,5),
    (8,8),
    (8,8),
    (9,9),
    (6,6),
    (9,9),
    (8,8),
    (9,9),
    (7,7),
    (9,9),
    (9,9),
    (9,9),
    (9,9),
    (8,8),
    (9,9),
    (9,9),
    (8,9),
    (9,9),
    (9,9),
    (10,10),
    (9,9),
    (9,9),
    (8,9),
    (9,9),
    (9,9),
    (9,9),
    (9,9),
    (9,9),
    (9,9),
    (9,9),
    (9,9),
    (8,8),
    (9,9),
    (9,9),
    (9,9),
    (9,9),
    (9,9),
    (6,6),
    (5,5),
    (7,7),
    (8,8),
    (6,6),
    (6,6),
    (5,5),
    (8,8),
    (6,6),
    (5,5)
)

$startRow = 2
for ($k = 0; $k -lt $data.Count; $k++) {
    $r = $startRow + $k
    $pair = $data[$k]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
